# Applies the "3rd time using it" update: refreshes the Reports sheet
# with a new 36-row dataset (rows 2-37) of Firm / Time / Lawyers Registered.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Schoenherr", "47s", "2"),
  @("Cassidy Levy Kent", "44s", "0"),
  @("Magnusson Law", "12s", "0"),
  @("Reliance Corporate Advisors", "22s", "1"),
  @("Fox And Mandal", "11s", "0"),
  @("Spruson And Ferguson", "35s", "3"),
  @("Winston And Strawn", "07s", "2"),
  @("White And Case", "02min 56s", "3"),
  @("Walkers", "25s", "3"),
  @("McMillan", "06min 52s", "0"),
  @("Meitar Law", "16s", "1"),
  @("Tiruchelvam Associates", "06s", "1"),
  @("HFW", "03min 47s", "0"),
  @("Crowell And Moring", "31s", "3"),
  @("Dittmar And Indrenius", "22s", "0"),
  @("Tuca Zbarcea", "05s", "1"),
  @("MZM Legal", "11s", "1"),
  @("Titov", "09s", "1"),
  @("Taylor Wessing", "42s", "3"),
  @("Ropes And Gray", "17s", "3"),
  @("Pedersoli", "01min 04s", "0"),
  @("LEX Logmannsstofa", "04s", "1"),
  @("Finreg 360", "06s", "1"),
  @("ShinAndKim", "17s", "1"),
  @("NPP Legal", "05s", "1"),
  @("Spencer West", "26s", "3"),
  @("Watson Farley And Williams", "33s", "3"),
  @("Stikeman Elliott", "01min 19s", "0"),
  @("Panetta Consulting Group", "08s", "1"),
  @("Njord", "06s", "2"),
  @("White And Case", "06min 17s", "0"),
  @("VB Advocates", "05s", "1"),
  @("SEUM", "23s", "1"),
  @("Malley And Co", "02s", "1"),
  @("Matheson", "07s", "1"),
  @("Wolf Theiss", "17s", "3")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $rowData = $data[$i]
  $ws.Cells.Item($r, 1).Value = $rowData[0]
  $ws.Cells.Item($r, 2).Value = $rowData[1]
  # leading apostrophe forces text storage, matching the source data's
  # text-typed "Lawyers Registered" numbers (e.g. "1", "0", "12") instead
  # of auto-converting them to numeric cells.
  $ws.Cells.Item($r, 3).Value = "'" + $rowData[2]
}
